$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the new "Save" column header in H1, copying the formatting used by the
# other header cells (e.g. G1) so the new column matches the existing style.
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("H1").Value = "Save"

# Populate the Save column with 0 for every data row (rows 2-11)
for ($r = 2; $r -le 11; $r++) {
    $ws.Cells.Item($r, 8).Value = 0
}
